$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit: Word relocated the hidden "_GoBack" bookmark (which marks
# the location of the user's last edit) from the end of the document
# to a brand-new, otherwise-empty paragraph right after the "Nuestro
# equipo" heading paragraph (i.e. the existing empty paragraph that
# follows it gets split in two, and the new second half carries the
# bookmark).
# ------------------------------------------------------------------

# Step 1: drop the existing "_GoBack" bookmark whose current position is
# at the very end of the document (it is about to be relocated).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 2: find the heading paragraph ("Nuestro equipo") and the empty
# paragraph that immediately follows it - that's the paragraph we need
# to split so the new half can hold the relocated bookmark.
$headingRange = $d.Content
$headingRange.Find.Execute("Nuestro equipo", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0) | Out-Null
$headingPara = $d.Paragraphs(1)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Nuestro equipo`r") {
        $headingPara = $p
        break
    }
}
$emptyPara = $headingPara.Next()

# Collapsed insertion point at the end of that empty paragraph (i.e.
# right before its paragraph mark) - this is where the split happens.
$splitPoint = $d.Range($emptyPara.Range.End, $emptyPara.Range.End)

# Step 3: insert a brand-new paragraph (formatted the same as its
# neighbours) that contains only the relocated bookmark - no run/text.
# A trailing bare "<w:p/>" is included only to force Word to treat the
# inserted fragment as a genuine paragraph break (it is removed again
# right after).
$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$splitPoint.InsertXML($xml)

# Step 4: InsertXML always needs a trailing paragraph mark to force the
# real split, which leaves one extra bare empty paragraph behind -
# remove it now that the split has happened.
$strayPara = $emptyPara.Next().Next()
$strayPara.Range.Delete()
